# "add 2d act camera" - update the default scene's (row 2, villageScene /
# PioneerNoob) camera offset position + rotation so the 2D act camera is
# used instead of the old 3D-ish offset.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column J = CamOffestPos, Column K = CamOffestRot (see header row 1)
$ws.Range("J2").Value = "0,4.2,5.5"
$ws.Range("K2").Value = "25,180"

# Leave the selection on the cell that was last touched, matching the
# saved workbook's cursor position.
$ws.Range("K2").Select()
